$d = $word.ActiveDocument

# The paragraph currently reads (across 3 runs, with spellcheck proofErr
# markers bracketing the middle run):
#   "My second " + "file" + " "
# It needs to become (2 runs, no proofErr markers):
#   "My second file " + "agduguagha"

# Step 1: rewrite the whole visible span (all three original runs) in one
# go. Because the replacement text differs from the original, the engine
# re-serialises this span as a single run and drops the now-enclosed
# w:proofErr spellStart/spellEnd markers along with it. "X" is a
# placeholder standing in for the new trailing word so step 1 and step 2
# remain independent, simple text swaps.
$whole = $d.Range(0, 15)
$whole.Text = "My second file X"

# Step 2: turn the placeholder into the real word "agduguagha".
$tail = $d.Range(15, 16)
$tail.Text = "agduguagha"

# Step 3: the text above is still one merged run (same rPr all the way
# through). The source paragraph wants the new word in its own run, so
# nudge the formatting of that tail range off of, and back to, its
# original value. That forces the engine to split it into its own run
# without leaving any visible formatting change behind.
$tailRange = $d.Range(15, 25)
$tailRange.Bold = 1
$tailRange.Bold = 0
